# Add "2022-Q3" sheet data and update the "总计" (summary) sheet.
#
# Layout before the edit:
#   1: 总计     2: 2022-Q2   3: 2022-Q1   4: 2021-Q4
#   5: 2021-Q3  6: 2021-Q2   7: 2021-Q1   8: 2020-Q4
#
# Layout after the edit:
#   1: 总计     2: 2022-Q3   3: 2022-Q2   4: 2022-Q1   5: 2021-Q4
#   6: 2021-Q3  7: 2021-Q2   8: 2021-Q1   9: 2020-Q4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the brand-new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $totalSheet)
$newSheet.Name = "2022-Q3"

# Touching UsedRange "wakes up" the freshly-created sheet so that the
# PasteSpecial calls below reliably keep their formatting.
$null = $newSheet.UsedRange

# The sheet that used to be "2022-Q2" is now pushed one slot further down
# (index 3); re-fetch it *after* the Add() above so the reference isn't
# stale, and reuse its header/data formatting for the new sheet.
$templateSheet = $wb.Worksheets.Item(3)
$templateSheet.Range("A1:H4").Copy()
$newSheet.Range("A1:H4").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B (fund code) and D:G (numeric-looking figures) are stored as
# plain text in this workbook (e.g. to keep leading zeros in fund codes) -
# mark them as Text before writing so Excel doesn't coerce them to numbers.
$newSheet.Range("B2:B4").NumberFormat = "@"
$newSheet.Range("D2:G4").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "159617"
$newSheet.Range("C2").Value = "华夏中证智选500价值稳健策略ETF"
$newSheet.Range("D2").Value = "2.93"
$newSheet.Range("E2").Value = "97.05"
$newSheet.Range("F2").Value = "1.50"
$newSheet.Range("G2").Value = "0.0440"
$newSheet.Range("H2").Value = 1

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "010154"
$newSheet.Range("C3").Value = "中加中证500指数增强C"
$newSheet.Range("D3").Value = "0.51"
$newSheet.Range("E3").Value = "94.15"
$newSheet.Range("F3").Value = "1.51"
$newSheet.Range("G3").Value = "0.0077"
$newSheet.Range("H3").Value = 8

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "010153"
$newSheet.Range("C4").Value = "中加中证500指数增强A"
$newSheet.Range("D4").Value = "0.44"
$newSheet.Range("E4").Value = "94.15"
$newSheet.Range("F4").Value = "1.51"
$newSheet.Range("G4").Value = "0.0066"
$newSheet.Range("H4").Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q3 and
#    renumber the leading index column.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Re-apply the existing formatting (bold + border on col A, plain on B:D)
# to the freshly inserted row by copying it from the row right below.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)     # xlPasteFormats
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.06

# Renumber the index column (A3:A9) sequentially: 1,2,3,4,5,6,7
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
